$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1581.1875
$ws.Range("I18").Value = 1831.6923
$ws.Range("J18").Value = 495.66666
$ws.Range("K18").Value = 1831.6923
$ws.Range("L18").Value = 495.66666
$ws.Range("M18").Value = -1547.6923
$ws.Range("N18").Value = -1063.66666
$ws.Range("H19").Value = 456.16666
$ws.Range("J19").Value = 414.16666
$ws.Range("L19").Value = 414.16666
$ws.Range("N19").Value = -764.16666
$ws.Range("H20").Value = 1088
$ws.Range("I20").Value = 1088
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1088
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -858
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("H35").Value = 1088
$ws.Range("I35").Value = 1088
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1088
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -709
$ws.Range("H48").Value = 1407.1578
$ws.Range("I48").Value = 871.58826
$ws.Range("K48").Value = 2614.76478
$ws.Range("M48").Value = -2322.76478
$ws.Range("H56").Value = 1407.1578
$ws.Range("I56").Value = 871.58826
$ws.Range("K56").Value = 2614.76478
$ws.Range("M56").Value = -2080.76478
$ws.Range("N20").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("N35").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30109.627
$ws.Range("I32").Value = 6583.184
$ws.Range("K32").Value = 6583.184
$ws.Range("M32").Value = -6296.184
$ws.Range("H54").Value = 29247.25
$ws.Range("J54").Value = 29247.25
$ws.Range("L54").Value = 29247.25
$ws.Range("N54").Value = -30785.25
$ws.Range("H74").Value = 1658.683
$ws.Range("I74").Value = 1446.3429
$ws.Range("K74").Value = 1446.3429
$ws.Range("M74").Value = -572.3429000000001
$ws.Range("H77").Value = 1658.683
$ws.Range("I77").Value = 1446.3429
$ws.Range("K77").Value = 7231.7145
$ws.Range("M77").Value = -2863.7145
$ws.Range("H122").Value = 1920.6177
$ws.Range("I122").Value = 1843.2858
$ws.Range("J122").Value = 2281.5
$ws.Range("K122").Value = 5529.857400000001
$ws.Range("L122").Value = 6844.5
$ws.Range("M122").Value = -3079.857400000001
$ws.Range("N122").Value = -11744.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 19646.834
$ws.Range("I128").Value = 19646.834
$ws.Range("K128").Value = 58940.50199999999
$ws.Range("M128").Value = -56450.50199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7500
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 7500
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
$ws.Range("H99").Value = 31462.812
$ws.Range("I99").Value = 33950.418
$ws.Range("K99").Value = 33950.418
$ws.Range("M99").Value = -32452.418
$ws.Range("H126").Value = 31462.812
$ws.Range("I126").Value = 33950.418
$ws.Range("K126").Value = 101851.254
$ws.Range("M126").Value = -99381.25399999999
$ws.Range("H132").Value = 2665.1724
$ws.Range("I132").Value = 2762.7
$ws.Range("K132").Value = 8288.099999999999
$ws.Range("M132").Value = -5758.099999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1649.6666
$ws.Range("I14").Value = 1649.6666
$ws.Range("K14").Value = 4948.9998
$ws.Range("M14").Value = -4775.9998
$ws.Range("H68").Value = 1599.375
$ws.Range("I68").Value = 1824.75
$ws.Range("J68").Value = 1374
$ws.Range("K68").Value = 5474.25
$ws.Range("L68").Value = 4122
$ws.Range("M68").Value = -4663.25
$ws.Range("N68").Value = -5744
$ws.Range("H71").Value = 1599.375
$ws.Range("I71").Value = 1824.75
$ws.Range("J71").Value = 1374
$ws.Range("K71").Value = 16422.75
$ws.Range("L71").Value = 12366
$ws.Range("M71").Value = -12366.75
$ws.Range("N71").Value = -20478
$ws.Range("H122").Value = 646.6429000000001
$ws.Range("J122").Value = 760.2222
$ws.Range("L122").Value = 6841.999800000001
$ws.Range("N122").Value = -11741.9998
$ws.Range("H138").Value = 791.6667
$ws.Range("I138").Value = 791.6667
$ws.Range("K138").Value = 2375.0001
$ws.Range("M138").Value = 2764.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 13900
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H36").Value = 820
$ws.Range("I36").Value = 820
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 820
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -335
$ws.Range("H43").Value = 25714.594
$ws.Range("J43").Value = 34779.914
$ws.Range("L43").Value = 34779.914
$ws.Range("N43").Value = -35081.914
$ws.Range("H132").Value = 3256.2454
$ws.Range("I132").Value = 2765.5454
$ws.Range("J132").Value = 5655.222
$ws.Range("K132").Value = 8296.636200000001
$ws.Range("L132").Value = 16965.666
$ws.Range("M132").Value = -5766.636200000001
$ws.Range("N132").Value = -22025.666
$ws.Range("N35").ClearContents()
$ws.Range("N36").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -187
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 300
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 300
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -68
$ws.Range("H37").Value = 300
$ws.Range("I37").Value = 300
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 300
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -193
$ws.Range("H100").Value = 24104.412
$ws.Range("I100").Value = 4966.4443
$ws.Range("J100").Value = 45634.625
$ws.Range("K100").Value = 4966.4443
$ws.Range("L100").Value = 45634.625
$ws.Range("M100").Value = -4425.4443
$ws.Range("N100").Value = -46716.625
$ws.Range("H122").Value = 7522.115
$ws.Range("J122").Value = 4768.4
$ws.Range("L122").Value = 14305.2
$ws.Range("N122").Value = -19205.2
$ws.Range("N4").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("N37").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 100000
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100462
$ws.Range("H51").Value = 12036.526
$ws.Range("J51").Value = 19763.182
$ws.Range("L51").Value = 19763.182
$ws.Range("N51").Value = -20783.182
$ws.Range("H62").Value = 61608.082
$ws.Range("I62").Value = 4100.6665
$ws.Range("J62").Value = 80777.22
$ws.Range("K62").Value = 4100.6665
$ws.Range("L62").Value = 80777.22
$ws.Range("M62").Value = -3476.6665
$ws.Range("N62").Value = -82025.22
$ws.Range("H65").Value = 61608.082
$ws.Range("I65").Value = 4100.6665
$ws.Range("J65").Value = 80777.22
$ws.Range("K65").Value = 20503.3325
$ws.Range("L65").Value = 403886.1
$ws.Range("M65").Value = -17383.3325
$ws.Range("N65").Value = -410126.1
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070
